$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 44.078976249067999
$ws.Range("B3").Value = 51.092697203883503
$ws.Range("B4").Value = 54.373902589237602
$ws.Range("B5").Value = 48.271624963585801
$ws.Range("B6").Value = 46.651265992614398
$ws.Range("B7").Value = 63.5080190545874
$ws.Range("B8").Value = 52.281367097161201
$ws.Range("B9").Value = 65.211545584619401
$ws.Range("B10").Value = 49.752055288256798
$ws.Range("B11").Value = 51.297537267271501
$ws.Range("B12").Value = 60.479717633469001
$ws.Range("B13").Value = 56.787098488378902
$ws.Range("B14").Value = 58.2017214476692
$ws.Range("B15").Value = 50.168146757817198
$ws.Range("B16").Value = 65.879618127352202
$ws.Range("B17").Value = 57.6138404682115
$ws.Range("B18").Value = 55.792555017753102
$ws.Range("B19").Value = 49.396279287752101
$ws.Range("B20").Value = 62.547170327540897
$ws.Range("B21").Value = 54.028989436674998
$ws.Range("B22").Value = 45.428433030854997
$ws.Range("B23").Value = 54.768350085038598
$ws.Range("B24").Value = 48.4073891916075
$ws.Range("B25").Value = 56.1644475536933
$ws.Range("B26").Value = 58.080732087846997
$ws.Range("B27").Value = 67.474923809769905
$ws.Range("B28").Value = 55.812893231662301
$ws.Range("B29").Value = 52.576520797828501
$ws.Range("B30").Value = 52.8140607390662
$ws.Range("B31").Value = 40.916123955462403
$ws.Range("B32").Value = 51.057367680950897
$ws.Range("B33").Value = 44.119733985371099
$ws.Range("B34").Value = 56.298681243666202
$ws.Range("B35").Value = 60.147762835105198
$ws.Range("B36").Value = 62.738564702866903
$ws.Range("B37").Value = 57.083664020299999
$ws.Range("B38").Value = 61.224693322117197
$ws.Range("B39").Value = 68.407517099340694
$ws.Range("B40").Value = 55.539008821503202
$ws.Range("B41").Value = 50.765420171254
$ws.Range("B42").Value = 46.403743894956698
$ws.Range("B43").Value = 48.214896516080401
$ws.Range("B44").Value = 38.801521577867099
$ws.Range("B45").Value = 42.476441402474997
$ws.Range("B46").Value = 67.070963512948396
$ws.Range("B47").Value = 58.6925956278806
$ws.Range("B48").Value = 62.956100730876301
$ws.Range("B49").Value = 66.3412933869715
$ws.Range("B50").Value = 54.5443455145248
$ws.Range("B51").Value = 55.371661756383297
$ws.Range("B52").Value = 49.223044252675201
$ws.Range("B53").Value = 47.263367371011
$ws.Range("B54").Value = 57.438808823446799
$ws.Range("B55").Value = 43.489130256925598
$ws.Range("B56").Value = 66.188051650142299
$ws.Range("B57").Value = 58.435485133667797
$ws.Range("B58").Value = 37.161129504195401
$ws.Range("B59").Value = 63.6969273197452
$ws.Range("B60").Value = 56.257447568449699
$ws.Range("B61").Value = 65.1800537091911
$ws.Range("B62").Value = 51.984063993124799
$ws.Range("B63").Value = 53.489477301848297
$ws.Range("B64").Value = 51.748004995355402
$ws.Range("B65").Value = 59.527971266005402
$ws.Range("B66").Value = 56.129470470847401
$ws.Range("B67").Value = 62.957150244474697
$ws.Range("B68").Value = 48.432257613953297
$ws.Range("B69").Value = 65.182750418630704
$ws.Range("B70").Value = 54.582238020310001
$ws.Range("B71").Value = 51.1250002212821
$ws.Range("B72").Value = 63.261863760880203
$ws.Range("B73").Value = 46.331298979215198

$ws.Range("F9").Select()
